$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the scenario description (column D) for each test row: replace the
# long step-by-step instructions with the short "<Action> Setup Mapping
# Jurnal Investasi" wording.
$ws.Range("D2").Value = "Tambah Setup Mapping Jurnal Investasi"
$ws.Range("D3").Value = "View Setup Mapping Jurnal Investasi"
$ws.Range("D4").Value = "Ubah Setup Mapping Jurnal Investasi"
$ws.Range("D5").Value = "Hapus Setup Mapping Jurnal Investasi"

# The shorter text no longer needs the same wrapped row height as before;
# let Excel recompute it, then pin rows 2, 3 and 5 to their new (smaller)
# height. Row 4's text now fits on a single line, so it keeps the default
# sheet row height with no explicit override.
$ws.Rows("2:5").AutoFit()
$ws.Rows("2").RowHeight = 30
$ws.Rows("3").RowHeight = 30
$ws.Rows("5").RowHeight = 30

# Leave the selection on D5, matching where the editor finished working.
$ws.Range("D5").Select() | Out-Null
